$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1015.8077
$ws.Range("J17").Value = 1083.0869
$ws.Range("L17").Value = 3249.2607
$ws.Range("N17").Value = -3585.2607
# Row 28
$ws.Range("H28").Value = 726.3125
$ws.Range("I28").Value = 513.5
$ws.Range("J28").Value = 1081
$ws.Range("K28").Value = 513.5
$ws.Range("L28").Value = 1081
$ws.Range("M28").Value = -28.5
$ws.Range("N28").Value = -2051
# Row 43
$ws.Range("H43").Value = 690.2
$ws.Range("I43").Value = 637.625
$ws.Range("J43").Value = 714.94116
$ws.Range("K43").Value = 637.625
$ws.Range("L43").Value = 714.94116
$ws.Range("M43").Value = -568.625
$ws.Range("N43").Value = -852.94116
# Row 75
$ws.Range("H75").Value = 27750
$ws.Range("J75").Value = 27750
$ws.Range("L75").Value = 27750
$ws.Range("N75").Value = -29622
# Row 78
$ws.Range("H78").Value = 27750
$ws.Range("J78").Value = 27750
$ws.Range("L78").Value = 83250
$ws.Range("N78").Value = -92610
# Row 86
$ws.Range("H86").Value = 56084.453
$ws.Range("I86").Value = 93646.38
$ws.Range("J86").Value = 1828.3334
$ws.Range("K86").Value = 93646.38
$ws.Range("L86").Value = 1828.3334
$ws.Range("M86").Value = -92523.38
$ws.Range("N86").Value = -4074.3334
# Row 89
$ws.Range("H89").Value = 56084.453
$ws.Range("I89").Value = 93646.38
$ws.Range("J89").Value = 1828.3334
$ws.Range("K89").Value = 468231.9
$ws.Range("L89").Value = 9141.666999999999
$ws.Range("M89").Value = -462615.9
$ws.Range("N89").Value = -20373.667
# Row 107
$ws.Range("H107").Value = 1060.6666
$ws.Range("I107").Value = 1041
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1041
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 879
$ws.Range("N107").Value = -4940
# Row 111
$ws.Range("H111").Value = 1065.9445
$ws.Range("I111").Value = 1095.8667
$ws.Range("J111").Value = 916.3333
$ws.Range("K111").Value = 3287.6001
$ws.Range("L111").Value = 2748.9999
$ws.Range("M111").Value = -220.6001000000001
$ws.Range("N111").Value = -8882.999899999999
# Row 113
$ws.Range("H113").Value = 3381.8
$ws.Range("I113").Value = 2740
$ws.Range("J113").Value = 4344.5
$ws.Range("K113").Value = 2740
$ws.Range("L113").Value = 4344.5
$ws.Range("M113").Value = 514
$ws.Range("N113").Value = -10852.5
# Row 114
$ws.Range("H114").Value = 48995
$ws.Range("J114").Value = 48995
$ws.Range("L114").Value = 48995
$ws.Range("N114").Value = -57673
# Row 116
$ws.Range("H116").Value = 3331
$ws.Range("I116").Value = 2996.5
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 2996.5
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 445.5
$ws.Range("N116").Value = -10884
# Row 120
$ws.Range("H120").Value = 38400
$ws.Range("J120").Value = 38400
$ws.Range("L120").Value = 38400
$ws.Range("N120").Value = -48076
# Row 129
$ws.Range("H129").Value = 1629.8096
$ws.Range("J129").Value = 2107.5862
$ws.Range("L129").Value = 6322.758600000001
$ws.Range("N129").Value = -16322.7586

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 9500
$ws.Range("J3").Value = 9500
$ws.Range("L3").Value = 9500
$ws.Range("N3").Value = -9730
# Row 22
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 10500
$ws.Range("I22").Value = 10500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -10201
# Row 32
$ws.Range("H32").Value = 5193.61
$ws.Range("I32").Value = 4312.9893
$ws.Range("J32").Value = 18990
$ws.Range("K32").Value = 4312.9893
$ws.Range("L32").Value = 18990
$ws.Range("M32").Value = -4025.9893
$ws.Range("N32").Value = -19564
# Row 41
$ws.Range("N41").ClearContents()
$ws.Range("H41").Value = 3185.3333
$ws.Range("I41").Value = 3185.3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3185.3333
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2771.3333
# Row 45
$ws.Range("H45").Value = 1084.7693
$ws.Range("I45").Value = 1091.8334
$ws.Range("K45").Value = 1091.8334
$ws.Range("M45").Value = -714.8334
# Row 56
$ws.Range("H56").Value = 10166.667
$ws.Range("J56").Value = 10166.667
$ws.Range("L56").Value = 10166.667
$ws.Range("N56").Value = -11650.667
# Row 61
$ws.Range("H61").Value = 1762.6842
$ws.Range("I61").Value = 1606.5555
$ws.Range("K61").Value = 1606.5555
$ws.Range("M61").Value = -1394.5555
# Row 132
$ws.Range("H132").Value = 1899201.1
$ws.Range("I132").Value = 6249.607
$ws.Range("K132").Value = 18748.821
$ws.Range("M132").Value = -16218.821
# Row 136
$ws.Range("H136").Value = 1762.6842
$ws.Range("I136").Value = 1606.5555
$ws.Range("K136").Value = 4819.666499999999
$ws.Range("M136").Value = -2269.666499999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 10000
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5280
# Row 94
$ws.Range("H94").Value = 1643.0714
$ws.Range("I94").Value = 1599.8182
$ws.Range("J94").Value = 1801.6666
$ws.Range("K94").Value = 1599.8182
$ws.Range("L94").Value = 1801.6666
$ws.Range("M94").Value = -1148.8182
$ws.Range("N94").Value = -2703.6666
# Row 107
$ws.Range("H107").Value = 1782.9615
$ws.Range("I107").Value = 1764.875
$ws.Range("K107").Value = 1764.875
$ws.Range("M107").Value = 155.125
# Row 130
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 15381.8
$ws.Range("J43").Value = 15381.8
$ws.Range("L43").Value = 15381.8
$ws.Range("N43").Value = -15749.8
# Row 94
$ws.Range("H94").Value = 1669.7142
$ws.Range("J94").Value = 1648.3334
$ws.Range("L94").Value = 1648.3334
$ws.Range("N94").Value = -2550.3334
# Row 101
$ws.Range("H101").Value = 15381.8
$ws.Range("J101").Value = 15381.8
$ws.Range("L101").Value = 15381.8
$ws.Range("N101").Value = -21871.8
# Row 105
$ws.Range("H105").Value = 1726.9231
$ws.Range("I105").Value = 1789.1666
$ws.Range("J105").Value = 980
$ws.Range("K105").Value = 1789.1666
$ws.Range("L105").Value = 980
$ws.Range("M105").Value = -42.16660000000002
$ws.Range("N105").Value = -4474
# Row 134
$ws.Range("H134").Value = 3156.4707
$ws.Range("I134").Value = 1636
$ws.Range("J134").Value = 3790
$ws.Range("K134").Value = 4908
$ws.Range("L134").Value = 11370
$ws.Range("M134").Value = -2373
$ws.Range("N134").Value = -16440
# Row 138
$ws.Range("H138").Value = 40796.25
$ws.Range("J138").Value = 42232.273
$ws.Range("L138").Value = 42232.273
$ws.Range("N138").Value = -52512.273
# Row 139
$ws.Range("H139").Value = 32960
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 34933.332
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 34933.332
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -45213.332

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 3371.4285
$ws.Range("J34").Value = 3371.4285
$ws.Range("L34").Value = 10114.2855
$ws.Range("N34").Value = -10282.2855
# Row 113
$ws.Range("H113").Value = 891.2857
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 914.8333
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 2744.4999
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -7084.4999
# Row 122
$ws.Range("H122").Value = 1720.8616
$ws.Range("J122").Value = 2048.0981
$ws.Range("L122").Value = 18432.8829
$ws.Range("N122").Value = -23332.8829
# Row 131
$ws.Range("H131").Value = 884.9722
$ws.Range("I131").Value = 450.66666
$ws.Range("J131").Value = 1195.1904
$ws.Range("K131").Value = 1351.99998
$ws.Range("L131").Value = 3585.5712
$ws.Range("M131").Value = 3688.00002
$ws.Range("N131").Value = -13665.5712
# Row 132
$ws.Range("H132").Value = 1178.0938
$ws.Range("J132").Value = 1237.2106
$ws.Range("L132").Value = 11134.8954
$ws.Range("N132").Value = -16194.8954

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 1692.5834
$ws.Range("I113").Value = 1692.5834
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1692.5834
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 477.4166
# Row 122
$ws.Range("H122").Value = 2465.9412
$ws.Range("I122").Value = 2744.3572
$ws.Range("J122").Value = 1166.6666
$ws.Range("K122").Value = 8233.071599999999
$ws.Range("L122").Value = 3499.9998
$ws.Range("M122").Value = -5783.071599999999
$ws.Range("N122").Value = -8399.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2519.8
$ws.Range("I61").Value = 2322.2727
$ws.Range("J61").Value = 3968.3333
$ws.Range("K61").Value = 2322.2727
$ws.Range("L61").Value = 3968.3333
$ws.Range("M61").Value = -2120.2727
$ws.Range("N61").Value = -4372.3333
# Row 113
$ws.Range("H113").Value = 2519.8
$ws.Range("I113").Value = 2322.2727
$ws.Range("J113").Value = 3968.3333
$ws.Range("K113").Value = 2322.2727
$ws.Range("L113").Value = 3968.3333
$ws.Range("M113").Value = -152.2727
$ws.Range("N113").Value = -8308.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 8101.091
$ws.Range("J69").Value = 8101.091
$ws.Range("L69").Value = 8101.091
$ws.Range("N69").Value = -9599.091
# Row 72
$ws.Range("H72").Value = 8101.091
$ws.Range("J72").Value = 8101.091
$ws.Range("L72").Value = 24303.273
$ws.Range("N72").Value = -31791.273
# Row 121
$ws.Range("H121").Value = 27616.666
$ws.Range("J121").Value = 27616.666
$ws.Range("L121").Value = 27616.666
$ws.Range("N121").Value = -31110.666
